# ACPerformanceManager: Mission profile - update Costs.xlsx DOC values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - Depreciation
$ws.Range("B2").Value = 2684.6999999999985
$ws.Range("C2").Value = 1917.6428571428569
$ws.Range("D2").Value = 6.711749999999996
$ws.Range("E2").Value = 5.162884615384613

# Row 3 - Interest
$ws.Range("B3").Value = 2577.3119999999985
$ws.Range("C3").Value = 1840.9371428571426
$ws.Range("D3").Value = 6.443279999999996
$ws.Range("E3").Value = 4.956369230769228

# Row 4 - Insurance
$ws.Range("B4").Value = 207.7333333333333
$ws.Range("C4").Value = 148.3809523809524
$ws.Range("D4").Value = 0.5193333333333332
$ws.Range("E4").Value = 0.39948717948717943

# Row 5 - DOC Capital
$ws.Range("B5").Value = 5469.74533333333
$ws.Range("C5").Value = 3906.9609523809513
$ws.Range("D5").Value = 13.674363333333325
$ws.Range("E5").Value = 10.518741025641019

# Row 7 - Cockpit Crew
$ws.Range("B7").Value = 1007.9999999999997
$ws.Range("D7").Value = 2.519999999999999
$ws.Range("E7").Value = 1.9384615384615378

# Row 8 - Cabin Crew
$ws.Range("B8").Value = 377.99999999999983
$ws.Range("D8").Value = 0.9449999999999996
$ws.Range("E8").Value = 0.7269230769230766

# Row 9 - DOC Crew
$ws.Range("B9").Value = 1385.9999999999995
$ws.Range("D9").Value = 3.464999999999999
$ws.Range("E9").Value = 2.665384615384614

# Row 11 - DOC Fuel
$ws.Range("B11").Value = 672.9592695723795
$ws.Range("C11").Value = 480.68519255169986
$ws.Range("D11").Value = 1.6823981739309488
$ws.Range("E11").Value = 1.2941524414853451

# Row 13 - Landing charges
$ws.Range("C13").Value = 363.16242857142856

# Row 14 - Navigation charges
$ws.Range("C14").Value = 302.082184000666

# Row 15 - Ground handling charges
$ws.Range("C15").Value = 947.1428571428572

# Row 16 - DOC Charges
$ws.Range("C16").Value = 1612.3874697149513

# Row 18 - Airframe Maintenance Charges
$ws.Range("B18").Value = 1025.9848466230153
$ws.Range("C18").Value = 732.8463190164398
$ws.Range("D18").Value = 2.5649621165575383
$ws.Range("E18").Value = 1.973047781967337

# Row 19 - Engine Maintenance Charges
$ws.Range("B19").Value = 1281.7943920165299
$ws.Range("C19").Value = 915.5674228689502
$ws.Range("D19").Value = 3.2044859800413246
$ws.Range("E19").Value = 2.4649892154164035

# Row 20 - DOC Maintenance
$ws.Range("B20").Value = 2420.485693246031
$ws.Range("C20").Value = 1728.9183523185943
$ws.Range("D20").Value = 6.051214233115077
$ws.Range("E20").Value = 4.654780179319291

# Row 23 - Total DOC
$ws.Range("B23").Value = 12206.532753752672
$ws.Range("C23").Value = 8718.951966966197
$ws.Range("D23").Value = 30.516331884381678
$ws.Range("E23").Value = 23.47410144952437

# Row 25 - Cash DOC
$ws.Range("B25").Value = 6736.787420419341
$ws.Range("C25").Value = 4811.9910145852455
$ws.Range("D25").Value = 16.841968551048353
$ws.Range("E25").Value = 12.955360423883349
